$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The B2:B101 values need to be rotated "up" by 6 rows (with wraparound):
# new B(r) = old B(r+6) for r in 2..95
# new B(r) = old B(r+6-100) for r in 96..101 (wraps to the original top values)

$rng = $ws.Range("B2:B101")
$orig = $rng.Value2

$n = 100
$shift = 6
$new = New-Object 'object[,]' $n,1
for ($i = 0; $i -lt $n; $i++) {
    $srcIndex = ($i + $shift) % $n
    $new[$i,0] = $orig[$srcIndex + 1, 1]
}

$rng.Value2 = $new

# Update the view: scroll so that A66 is the top-left visible cell,
# and select B101 as the active cell.
$ws.Range("B101").Select()
$excel.ActiveWindow.ScrollRow = 66
$excel.ActiveWindow.ScrollColumn = 1
